# Append the new "Billing foundation implemented" update block to the end of the
# requirement-status document, right after the last existing paragraph and before
# the closing section properties.
$d = $word.ActiveDocument

$newLines = @(
    "",
    "Update: 2026-02-21 (Billing foundation implemented)",
    "- Added billing domain models and schema:",
    "  - BillingPlans, TenantSubscriptions, TenantUsages, BillingInvoices",
    "- Added plan/subscription/usage/invoice seed data for demo tenants.",
    "- Added platform billing plan CRUD API:",
    "  - GET/POST/PUT/DELETE /api/platform/billing/plans",
    "- Added user billing APIs:",
    "  - GET /api/billing/plans",
    "  - GET /api/billing/current-plan",
    "  - GET /api/billing/usage",
    "  - GET /api/billing/invoices",
    "  - POST /api/billing/change-plan",
    "  - POST /api/billing/cancel",
    "- Added public plans API for landing sync:",
    "  - GET /api/public/plans",
    "- Wired billing limit enforcement service and enabled team member limit check on invite.",
    "- Wired frontend Billing page to real backend billing APIs (plan, usage, invoices, change/cancel).",
    "- Added Platform Settings Billing Plans tab with create/edit/archive + limits/features configuration.",
    "- Landing page pricing now fetches plans from backend public endpoint with fallback."
)

$insertionPoint = $d.Content
$insertionPoint.Collapse(0) | Out-Null

foreach ($line in $newLines) {
    $insertionPoint.InsertAfter("`r" + $line)
    $insertionPoint.Collapse(0) | Out-Null
}

Write-Output ("Paragraphs now in document: " + $d.Paragraphs.Count)
